$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "ProductEmail"

$ws2.Range("B1").Value = "To"
$ws2.Range("C1").Value = "From"
$ws2.Range("D1").Value = "Subject"
$ws2.Range("E1").Value = "Message"

$ws2.Range("B2").Value = "rais@softway.com"
$ws2.Range("C2").Value = "raees@softway.com"
$ws2.Range("D2").Value = "Test"
$ws2.Range("E2").Value = "This is a test message..."

$ws2.Range("B3").Value = "raees@softway.com"
$ws2.Range("C3").Value = "rais@softway.com"
$ws2.Range("D3").Value = "Test"
$ws2.Range("E3").Value = "This is a test message2..."

$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:rais@softway.com")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:raees@softway.com")
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:raees@softway.com")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "mailto:rais@softway.com")

Write-Host "done"
